$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table spans A2:AY13 (12 data rows, 51 columns, header on row 1).
# This edit rotates the rows: the last data row (row 13) moves to the top
# of the data block (row 2), and the previous rows 2-12 shift down to 3-13.

$rows = 12
$cols = 51

# Force every non-numeric / non-boolean column to Text format BEFORE writing,
# so that the COM layer doesn't "smart type" strings that look like numbers
# or dates (e.g. "1", "2021-10-06") into real numbers/dates on write-back.
$textCols = @("C","D","F","G","H","I","J","K","L","M","N","O","P","T","U","V","W","X","Y","Z","AA","AB","AC","AF","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY")
foreach ($colName in $textCols) {
    $ws.Range(($colName + "2:" + $colName + "13")).NumberFormat = "@"
}

$srcRange = $ws.Range("A2:AY13")
$v = $srcRange.Value()

$new = New-Object 'object[,]' $rows,$cols

# New row 1 (sheet row 2) <- old row 12 (sheet row 13)
for ($c = 1; $c -le $cols; $c++) {
    $new[0, $c-1] = $v[12, $c]
}

# New rows 2..12 (sheet rows 3..13) <- old rows 1..11 (sheet rows 2..12)
for ($r = 1; $r -le 11; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $new[$r, $c-1] = $v[$r, $c]
    }
}

$srcRange.Value = $new
